# Update NATMI LR-pair TPM output table with refreshed values (Uts2-Uts2r)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Uts2"
$ws.Range("C2").Value = "Uts2r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.126357
$ws.Range("H2").Value = 0.379071
$ws.Range("I2").Value = 0.1226892855998496
$ws.Range("J2").Value = 0.1226892855998496
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8480913333333334
$ws.Range("N2").Value = 2.544274
$ws.Range("O2").Value = 0.8068784050391488
$ws.Range("P2").Value = 0.8068784050391488
$ws.Range("Q2").Value = 0.107162276606
$ws.Range("R2").Value = 0.964460489454
$ws.Range("S2").Value = 0.09899533508019923
$ws.Range("T2").Value = 0.09899533508019923

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Uts2"
$ws.Range("C3").Value = "Uts2r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.126357
$ws.Range("H3").Value = 0.379071
$ws.Range("I3").Value = 0.1226892855998496
$ws.Range("J3").Value = 0.1226892855998496
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1543526666666667
$ws.Range("N3").Value = 0.463058
$ws.Range("O3").Value = 0.1468519115789487
$ws.Range("P3").Value = 0.1468519115789487
$ws.Range("Q3").Value = 0.019503539902
$ws.Range("R3").Value = 0.175531859118
$ws.Range("S3").Value = 0.01801715612059349
$ws.Range("T3").Value = 0.01801715612059349

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Uts2"
$ws.Range("C4").Value = "Uts2r"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.126357
$ws.Range("H4").Value = 0.379071
$ws.Range("I4").Value = 0.1226892855998496
$ws.Range("J4").Value = 0.1226892855998496
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.048633
$ws.Range("N4").Value = 0.145899
$ws.Range("O4").Value = 0.04626968338190256
$ws.Range("P4").Value = 0.04626968338190256
$ws.Range("Q4").Value = 0.006145119981
$ws.Range("R4").Value = 0.055306079829
$ws.Range("S4").Value = 0.005676794399056857
$ws.Range("T4").Value = 0.005676794399056857

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Uts2"
$ws.Range("C5").Value = "Uts2r"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.324908
$ws.Range("H5").Value = 0.9747239999999999
$ws.Range("I5").Value = 0.3154770246656372
$ws.Range("J5").Value = 0.3154770246656372
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8480913333333334
$ws.Range("N5").Value = 2.544274
$ws.Range("O5").Value = 0.8068784050391488
$ws.Range("P5").Value = 0.8068784050391488
$ws.Range("Q5").Value = 0.2755516589306667
$ws.Range("R5").Value = 2.479964930376
$ws.Range("S5").Value = 0.2545515984887056
$ws.Range("T5").Value = 0.2545515984887056

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Uts2"
$ws.Range("C6").Value = "Uts2r"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.324908
$ws.Range("H6").Value = 0.9747239999999999
$ws.Range("I6").Value = 0.3154770246656372
$ws.Range("J6").Value = 0.3154770246656372
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1543526666666667
$ws.Range("N6").Value = 0.463058
$ws.Range("O6").Value = 0.1468519115789487
$ws.Range("P6").Value = 0.1468519115789487
$ws.Range("Q6").Value = 0.05015041622133333
$ws.Range("R6").Value = 0.451353745992
$ws.Range("S6").Value = 0.04632840413138797
$ws.Range("T6").Value = 0.04632840413138797

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Uts2"
$ws.Range("C7").Value = "Uts2r"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.324908
$ws.Range("H7").Value = 0.9747239999999999
$ws.Range("I7").Value = 0.3154770246656372
$ws.Range("J7").Value = 0.3154770246656372
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.048633
$ws.Range("N7").Value = 0.145899
$ws.Range("O7").Value = 0.04626968338190256
$ws.Range("P7").Value = 0.04626968338190256
$ws.Range("Q7").Value = 0.015801250764
$ws.Range("R7").Value = 0.142211256876
$ws.Range("S7").Value = 0.0145970220455437
$ws.Range("T7").Value = 0.0145970220455437

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Uts2"
$ws.Range("C8").Value = "Uts2r"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5786293333333333
$ws.Range("H8").Value = 1.735888
$ws.Range("I8").Value = 0.5618336897345133
$ws.Range("J8").Value = 0.5618336897345133
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8480913333333334
$ws.Range("N8").Value = 2.544274
$ws.Range("O8").Value = 0.8068784050391488
$ws.Range("P8").Value = 0.8068784050391488
$ws.Range("Q8").Value = 0.4907305228124445
$ws.Range("R8").Value = 4.416574705312001
$ws.Range("S8").Value = 0.4533314714702441
$ws.Range("T8").Value = 0.4533314714702441

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Uts2"
$ws.Range("C9").Value = "Uts2r"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5786293333333333
$ws.Range("H9").Value = 1.735888
$ws.Range("I9").Value = 0.5618336897345133
$ws.Range("J9").Value = 0.5618336897345133
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1543526666666667
$ws.Range("N9").Value = 0.463058
$ws.Range("O9").Value = 0.1468519115789487
$ws.Range("P9").Value = 0.1468519115789487
$ws.Range("Q9").Value = 0.08931298061155556
$ws.Range("R9").Value = 0.803816825504
$ws.Range("S9").Value = 0.08250635132696725
$ws.Range("T9").Value = 0.08250635132696725

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Uts2"
$ws.Range("C10").Value = "Uts2r"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5786293333333333
$ws.Range("H10").Value = 1.735888
$ws.Range("I10").Value = 0.5618336897345133
$ws.Range("J10").Value = 0.5618336897345133
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.048633
$ws.Range("N10").Value = 0.145899
$ws.Range("O10").Value = 0.04626968338190256
$ws.Range("P10").Value = 0.04626968338190256
$ws.Range("Q10").Value = 0.028140480368
$ws.Range("R10").Value = 0.253264323312
$ws.Range("S10").Value = 0.02599586693730201
$ws.Range("T10").Value = 0.02599586693730201

